$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy format from an existing header cell (reuses same style)
# then set the text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF)
$data = @(
    @(5, 7),
    @(5, 7),
    @(5, 7),
    @(6, 6),
    @(7, 8),
    @(4, 5),
    @(7, 7),
    @(7, 8),
    @(8, 9),
    @(9, 9),
    @(5, 7),
    @(6, 8),
    @(10, 10),
    @(6, 7),
    @(8, 8),
    @(8, 8)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
